$wb = $excel.ActiveWorkbook

# --- Reorder sheets -------------------------------------------------------
# Target tab order: Germany, Belgium, Sheet2, Czech, Slot Cards 215 Panel
# Move "Slot Cards 215 Panel" to the very end so "Sheet2" shifts up to slot 3.
$slotCards = $wb.Worksheets.Item("Slot Cards 215 Panel")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$slotCards.Move($null, $lastSheet)

# --- Add the new "Czech" sheet --------------------------------------------
# Copy "Germany" (same legend/layout) and drop it right after "Sheet2",
# i.e. right before "Slot Cards 215 Panel" which is now last.
$germany = $wb.Worksheets.Item("Germany")
$sheet2 = $wb.Worksheets.Item("Sheet2")
$germany.Copy($null, $sheet2)

$czech = $wb.Worksheets.Item("Sheet2").Next
$czech.Name = "Czech"

# Sheet-specific values.
$czech.Range("B2").Value = "Czech Market"
$czech.Range("B4").Value = "NGC-3477/T1734"

# Column widths specific to the Czech sheet.
$czech.Columns.Item(2).ColumnWidth = 33.333333
$czech.Columns.Item(3).ColumnWidth = 14.666666
$czech.Columns.Item(4).ColumnWidth = 21.333333

# Selection state + make it the active/selected tab.
$czech.Range("A8").Select()
$czech.Activate()

# --- Fix up selections / tabSelected on the other sheets -------------------
$germanySheet = $wb.Worksheets.Item("Germany")
$germanySheet.Range("A8").Select()

$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Range("A8").Select()

$sheet2Again = $wb.Worksheets.Item("Sheet2")
$sheet2Again.Range("A1:D10").Select()

# Re-activate Czech last so it ends up as the active tab.
$czech = $wb.Worksheets.Item("Czech")
$czech.Activate()
$czech.Range("A8").Select()
